# Update sea fuel mix (init_fuel_mix sheet): split MGO/HFO shares 100/0 -> 55/45
$wb = $excel.ActiveWorkbook

$wsFuelMix = $wb.Worksheets.Item("init_fuel_mix")
$wsComments = $wb.Worksheets.Item("Comments")

# Sea / MGO row (row 9) and Sea / HFO row (row 10)
$wsFuelMix.Range("D9").Value = 55
$wsFuelMix.Range("D10").Value = 45

# Document the new source + methodology on the Comments sheet
$wsComments.Range("A5").Value = "Other source for init mode mix: https://dokumen.tips/documents/teknisk-vurdering-av-skip-og-av-infrastruktur-for-forsyning-av-drivstoff-.html"
$wsComments.Range("A6").Value = "Split HFO/MGO based on figure 5-2 in source above"

# Refresh on-screen selections to match where the author left off
$wsComments.Activate() | Out-Null
$wsComments.Range("A7").Select() | Out-Null

$wsFuelMix.Activate() | Out-Null
$wsFuelMix.Range("J15").Select() | Out-Null
